$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated evaluation metrics after adding weight to household surplus
$data = @{
    2 = @(0.02, 0.29, 0.172621695626245, 0.5477959902107863, 1570764594.539834)
    3 = @(0.02, 0.29, 0.172621695626253, 0.5477959902108043, 1570764594.539834)
    4 = @(0.02, 0.29, 0.1726216956262516, 0.5477959902108045, 1570764594.539834)
    5 = @(0.02, 0.29, 0.172621695626253, 0.5477959902108043, 1570764594.539834)
    6 = @(0.02, 0.29, 0.1726216956262518, 0.5477959902108046, 1570764594.539834)
    7 = @(0.02, 0.29, 0.1726216956262532, 0.5477959902108044, 1570764594.539834)
    8 = @(0.02, 0.29, 0.1726216956262526, 0.5477959902108029, 1570764594.539834)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
}
